$wb = $excel.ActiveWorkbook

# --- Rename the "for_individual" translation row to "for_member" ---
# table_specific_translations: row 16 held the string_token/text pair for
# generating an entitlement for an individual; it now documents the member
# flow instead (translations individual -> member).
$tst = $wb.Worksheets.Item("table_specific_translations")
$tst.Range("A16").Value = "for_member"
$tst.Range("B16").Value = "Entitlement Generated for Member"

# --- Switch the active/selected sheet to table_specific_translations ---
[void]$tst.Activate()
[void]$tst.Range("B16").Select()
